$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from an existing header cell (A1) to the new header cell F1,
# then set its text to "Modelo".
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "Modelo"

# Update the slightly recomputed numeric values in B2 and D2.
$ws.Range("B2").Value = 0.02950307763024238
$ws.Range("D2").Value = 0.1292627146720762

# Add the model description in the new F2 cell.
$ws.Range("F2").Value = "Pipeline(steps=[('model', GradientBoostingRegressor(n_estimators=150))])"
